$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.595.28'
$ws.Range("E2").Value = '  -1.54%  '

$ws.Range("D3").Value = '1.589.16'
$ws.Range("E3").Value = '  -1.86%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = '''210.63'
$ws.Range("E5").Value = '  -1.72%  '

$ws.Range("E6").Value = '  -1.47%  '

$ws.Range("E7").Value = '  +0.13%  '

$ws.Range("E8").Value = '  -2.18%  '

$ws.Range("D9").Value = '''0.0615'
$ws.Range("E9").Value = '  -1.38%  '

$ws.Range("D10").Value = '''19.54'
$ws.Range("E10").Value = '  -3.35%  '

$ws.Range("E11").Value = '  -1.16%  '

$ws.Range("D12").Value = '1.812.78'
$ws.Range("E12").Value = '  -1.80%  '

$ws.Range("D13").Value = '1.589.98'
$ws.Range("E13").Value = '  -1.69%  '

$ws.Range("E14").Value = '  -2.44%  '

$ws.Range("D15").Value = '''0.522'
$ws.Range("E15").Value = '  -3.56%  '

$ws.Range("D16").Value = '''64.61'
$ws.Range("E16").Value = '  +0.47%  '

$ws.Range("D17").Value = '26.603.93'
$ws.Range("E17").Value = '  -1.43%  '

$ws.Range("D18").Value = '0.0₃0725'
$ws.Range("E18").Value = '  -2.04%  '

$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = '''1.00'
$ws.Range("E19").Value = '  +0.10%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '''208.17'
$ws.Range("E20").Value = '  -3.44%  '

$ws.Range("D21").Value = '''6.71'
$ws.Range("E21").Value = '  -2.23%  '

$ws.Range("D22").Value = '''4.24'
$ws.Range("E22").Value = '  -2.86%  '

$ws.Range("D23").Value = '''2.35'
$ws.Range("E23").Value = '  -2.42%  '

$ws.Range("D24").Value = '''8.84'
$ws.Range("E24").Value = '  -1.65%  '

$ws.Range("D25").Value = '''146.62'
$ws.Range("E25").Value = '  -0.48%  '

$ws.Range("E26").Value = '  +0.27%  '

$ws.Range("D27").Value = '''7.23'
$ws.Range("E27").Value = '  -0.39%  '

$ws.Range("E28").Value = '  -2.83%  '

$ws.Range("D29").Value = '''15.26'
$ws.Range("E29").Value = '  -1.69%  '

$ws.Range("D30").Value = '''0.0507'
$ws.Range("E30").Value = '  +0.76%  '

$ws.Range("E31").Value = '  -1.89%  '

$ws.Range("D32").Value = '''3.21'
$ws.Range("E32").Value = '  -3.69%  '

$ws.Range("D33").Value = '''0.664'
$ws.Range("E33").Value = '  +22.00%  '

$ws.Range("D34").Value = '''2.89'
$ws.Range("E34").Value = '  -2.74%  '

$ws.Range("D35").Value = '1.306.74'
$ws.Range("E35").Value = '  -2.44%  '

$ws.Range("D36").Value = '''2.43'
$ws.Range("E36").Value = '  -1.09%  '

$ws.Range("E37").Value = '  -4.85%  '

$ws.Range("E38").Value = '  -2.43%  '

$ws.Range("D39").Value = '''0.827'
$ws.Range("E39").Value = '  -2.08%  '

$ws.Range("E40").Value = '  +0.10%  '

$ws.Range("D41").Value = '''0.790'
$ws.Range("E41").Value = '  -1.38%  '

$ws.Range("D42").Value = '''5.36'
$ws.Range("E42").Value = '  +2.82%  '

$ws.Range("E43").Value = '  -2.89%  '

$ws.Range("D44").Value = '''62.75'
$ws.Range("E44").Value = '  -3.34%  '

$ws.Range("D45").Value = '1.725.04'
$ws.Range("E45").Value = '  -1.80%  '

$ws.Range("D46").Value = '''89.76'
$ws.Range("E46").Value = '  -0.62%  '

$ws.Range("D47").Value = '''1.60'
$ws.Range("E47").Value = '  -0.27%  '

$ws.Range("D48").Value = '''0.837'
$ws.Range("E48").Value = '  -1.27%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '''0.0978'
$ws.Range("E49").Value = '  -1.65%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.0502'
$ws.Range("E50").Value = '  -2.00%  '

$ws.Range("D51").Value = '''7.51'
$ws.Range("E51").Value = '  -0.56%  '
